$ErrorActionPreference = 'Stop'
$d = $word.ActiveDocument

# Locate the empty paragraph just before the final bookmark paragraph
$anchorIndex = $d.Paragraphs.Count - 1
$anchor = $d.Paragraphs.Item($anchorIndex)
$insertRange = $anchor.Range
$insertRange.Collapse(0)

$newParagraphsXml = @(
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>19/10/2018</w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> : </w:t></w:r><w:r><w:t>(</w:t></w:r><w:r><w:t>cours en autonomie</w:t></w:r><w:r><w:t>)</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Continuité de la création des interfaces utilisateurs</w:t></w:r><w:r><w:t xml:space="preserve"> (deuxième et troisième pages)</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Mise à jour de la planification</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Ajout des fichiers documentations dans le GitHub</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pBdr><w:bottom w:val="single" w:sz="6" w:space="1" w:color="auto"/></w:pBdr></w:pPr></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:ind w:left="765"/></w:pPr></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>06/11/2018 :</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Fin de la réalisation des pages interfaces du quizz</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Mise à jour de la planification</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Ajout des fichiers documentations dans le GitHub</w:t></w:r></w:p>'
)

foreach ($xml in $newParagraphsXml) {
    $newRange = $insertRange.InsertParagraphAfter()
    $insertRange.Collapse(0)
    $newParaIndex = $anchorIndex + 1
    $newPara = $d.Paragraphs.Item($newParaIndex)
    [void]$newPara.Range.InsertXML($xml)
    $anchorIndex = $newParaIndex
    $anchor = $d.Paragraphs.Item($anchorIndex)
    $insertRange = $anchor.Range
    $insertRange.Collapse(0)
}

# Append the final bullet run into the last paragraph (before the bookmark)
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertBefore('Commencement de la mise en place des questions')

Write-Host 'Edit applied.'
